# Add a new "blank_material" price sheet, positioned right after "shirka"
# (and therefore right before "interierka"), pre-populated with the same
# header/formatting style used on the neighbouring price sheets plus its
# own price data (purchase price / wholesale / retail formulas).

$wb = $excel.ActiveWorkbook

$shirka = $wb.Worksheets.Item("shirka")

# ---------------------------------------------------------------------
# 1. Create the new sheet right after "shirka" and rename it.
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Add($null, $shirka)
$ws.Name = "blank_material"

# Re-fetch after the insert shifted indices, so this points at the
# right sheet regardless of how the host tracks worksheet identity.
$finishka = $wb.Worksheets.Item("finishka")

# ---------------------------------------------------------------------
# 2. Pull over the cell formatting used by "shirka" so the new sheet
#    matches the existing look (header row styles, body styles, the
#    "wrap text" style used for the two long product names).
# ---------------------------------------------------------------------
$shirka.Range("A1:E1").Copy()
$ws.Range("A1:E1").PasteSpecial(-4122)

$shirka.Range("A2:E8").Copy()
$ws.Range("A2:E8").PasteSpecial(-4122)

$shirka.Range("A9:E9").Copy()
$ws.Range("A8:E9").PasteSpecial(-4122)

$shirka.Range("A60").Copy()
$ws.Range("A10").PasteSpecial(-4122)

$excel.CutCopyMode = 0

# ---------------------------------------------------------------------
# 3. Header row.
# ---------------------------------------------------------------------
$ws.Range("B1").Value = "Название"
$ws.Range("C1").Value = "Цена закупки"
$ws.Range("D1").Value = "Цена продажи РА"
$ws.Range("E1").Value = "Цена продажи Розница"

# ---------------------------------------------------------------------
# 4. Data rows: # / name / purchase price, with РА & retail formulas.
# ---------------------------------------------------------------------
$names = @(
    "Баннер 440 грамм ламинированный",
    "Баннер 510 грамм литой",
    "Баннер 510 грамм литой (Черный оборот)",
    "Баннер 340 грамм ламинированный",
    "Пленка самоклеющаяся",
    "Блюбек",
    "Баннерная сетка 370 грамм",
    "Бумага 150 грамм СytiLight "
)
$prices = @(100, 125, 140, 165, 230, 110, 250, 120)

for ($i = 0; $i -lt 8; $i++) {
    $row = 2 + $i
    $ws.Range("A$row").Value = $i + 1
    $ws.Range("B$row").Value = $names[$i]
    $ws.Range("C$row").Value = $prices[$i]
}

# Row 2 gets its own (non-shared) formulas ...
$ws.Range("D2").Formula = "=C2*2"
$ws.Range("E2").Formula = "=D2*1.3"

# ... rows 3-9 share one formula definition each, same as the source sheets.
$ws.Range("D3:D9").Formula = "=C3*2"
$ws.Range("E3:E9").Formula = "=D3*1.3"

# ---------------------------------------------------------------------
# 5. Column widths / row heights to match the sibling price sheets.
# ---------------------------------------------------------------------
$ws.Columns.Item(1).ColumnWidth = 1.7369791666666665
$ws.Columns.Item(2).ColumnWidth = 57.592447916666664
$ws.Columns.Item(3).ColumnWidth = 15.736979166666666
$ws.Columns.Item(4).ColumnWidth = 21.022135416666668
$ws.Columns.Item(5).ColumnWidth = 27.736979166666668

$ws.Range("A1").RowHeight = 15.75
$ws.Range("A2:A10").RowHeight = 15

# ---------------------------------------------------------------------
# 6. Selections / active sheet bookkeeping.
# ---------------------------------------------------------------------
$shirka.Activate()
$shirka.Range("A1:F10").Select()

$finishka.Activate()
$finishka.Range("D29").Select()

$ws.Activate()
$ws.Range("D12").Select()
